# Daily attendance processing - 2026-01-12 21:35:21
# Swap the order of names in the "Recorded By" column (G) wherever the
# cell contains exactly the two-part combination "System, dnasr281@gmail.com"
# or "admin@admin.com, dnasr281@gmail.com", turning it into
# "dnasr281@gmail.com, System" / "dnasr281@gmail.com, admin@admin.com"
# respectively. Other combinations (single names, or already starting with
# dnasr281@gmail.com, or including backup@backdoor.com) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G
    $value = $cell.Value2

    if ($value -ne $null -and $value -is [string] -and $value.Contains(",")) {
        $parts = $value.Split(",")
        if ($parts.Length -eq 2) {
            $first = $parts[0].Trim()
            $second = $parts[1].Trim()

            if (($first -eq "System" -and $second -eq "dnasr281@gmail.com") -or
                ($first -eq "admin@admin.com" -and $second -eq "dnasr281@gmail.com")) {
                $cell.Value2 = "$second, $first"
            }
        }
    }
}
